# major accuracy check update
#  - Rename the shared string "E7420" -> "E7420L" (column G, s2cDNAProtocol)
#  - Replace the H2:H29 "=FALSE()" formulas with literal boolean FALSE values
#    (same displayed/stored value of FALSE, but stored as a boolean literal
#    instead of a volatile formula).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the protocol code text from E7420 to E7420L wherever it appears.
$used = $ws.UsedRange
foreach ($cell in $used.Cells) {
    if ($cell.Value() -eq "E7420") {
        $cell.Value = "E7420L"
    }
}

# 2) Replace the FALSE() formulas in column H with literal boolean values.
for ($row = 2; $row -le 29; $row++) {
    $cell = $ws.Cells.Item($row, 8)
    $cell.Value = $false
}
